# "Show route as graph"
#
# The second route table (rows 13-22, Line/Station already filled in)
# gets its Time column (C) populated, mirroring the first table above
# it (rows 2-11) - this is the data the new route graph will be built
# from. The active selection also moves off the table (was C13) to
# F12, where the graph is going.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 3).Value = 36.3
$ws.Cells.Item(15, 3).Value = 13.2
$ws.Cells.Item(16, 3).Value = 18.1
$ws.Cells.Item(17, 3).Value = 15.7
$ws.Cells.Item(18, 3).Value = 9.4
$ws.Cells.Item(19, 3).Value = 40.4
$ws.Cells.Item(20, 3).Value = 13.8
$ws.Cells.Item(21, 3).Value = 11.1
$ws.Cells.Item(22, 3).Value = 38.3

$ws.Range("F12").Select()
